# Refactor code structure + fix pip package issue (#9)
#
# Updates the lab3 test-report workbook:
#  - Summary!A2/C2: total-test and failed-test counts drop (386/317 -> 350/281)
#  - All / Failed sheets: the old per-route "is missing route X" assertions
#    for the routing-table checks are replaced by a single
#    "have the wrong number of routes" assertion per device, and the
#    subsequent DNS/ping rows are renumbered upward to fill the gap; the
#    now-superfluous trailing rows are removed entirely.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet: update the Total/Failed test counts.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A2").Value = "350"
$summary.Range("C2").Value = "281"

# ---------------------------------------------------------------------
# 2) Shared replacement block used by both the "All" and "Failed" sheets:
#    20 rows of (Description, Reason) - the "Passed" column (B) stays
#    "False" throughout and is left untouched.
# ---------------------------------------------------------------------
$newRows = @(
    @('Checking the routing table of as1r1', 'The routing table of as1r1 have the wrong number of routes: 0, expected: 8'),
    @('Checking the routing table of as1r2', 'The routing table of as1r2 have the wrong number of routes: 0, expected: 8'),
    @('Checking the routing table of as2r1', 'The routing table of as2r1 have the wrong number of routes: 0, expected: 8'),
    @('Checking the routing table of as2r2', 'The routing table of as2r2 have the wrong number of routes: 0, expected: 8'),
    @('Checking the routing table of as3r1', 'The routing table of as3r1 have the wrong number of routes: 0, expected: 9'),
    @('Checking the routing table of root', 'The routing table of root have the wrong number of routes: 0, expected: 2'),
    @('Checking the routing table of net', 'The routing table of net have the wrong number of routes: 0, expected: 2'),
    @('Checking the routing table of pc', 'The routing table of pc have the wrong number of routes: 0, expected: 2'),
    @('Checking the routing table of local', 'The routing table of local have the wrong number of routes: 0, expected: 2'),
    @('Checking on `root` that `1.1.0.2` is the authority for domain `.`', 'named not started in the startup file of `root`'),
    @('Checking on `root` that `1.1.0.2` is the authority for domain `.`', 'named not started in the startup file of `root`'),
    @('Checking on `local` that `1.1.0.2` is the authority for domain `.`', 'named not started in the startup file of `local`'),
    @('Checking on `net` that `2.1.0.2` is the authority for domain `net`', 'named not started in the startup file of `net`'),
    @('Checking that `3.2.0.2` is the local name server for device `as1r1`', '`resolv.conf` file not found for device `as1r1`'),
    @('Checking that `3.2.0.2` is the local name server for device `as1r2`', '`resolv.conf` file not found for device `as1r2`'),
    @('Checking that `3.2.0.2` is the local name server for device `as2r1`', '`resolv.conf` file not found for device `as2r1`'),
    @('Checking that `3.2.0.2` is the local name server for device `as2r2`', '`resolv.conf` file not found for device `as2r2`'),
    @('Checking that `3.2.0.2` is the local name server for device `as3r1`', '`resolv.conf` file not found for device `as3r1`'),
    @('Checking that `3.2.0.2` is the local name server for device `as3r2`', '`resolv.conf` file not found for device `as3r2`'),
    @('Checking that `3.2.0.2` is the local name server for device `pc`', '`resolv.conf` file not found for device `pc`')
)

# ---------------------------------------------------------------------
# 3) Apply the block + trim trailing rows on both "All" and "Failed".
#    "All"    : rows 332-387 -> 332-351 (replace 20, drop 36)
#    "Failed" : rows 263-318 -> 263-282 (replace 20, drop 36)
# ---------------------------------------------------------------------
$sheetConfigs = @(
    @{ Name = "All"; StartRow = 332; OldLastRow = 387 },
    @{ Name = "Failed"; StartRow = 263; OldLastRow = 318 }
)

foreach ($cfg in $sheetConfigs) {
    $sheet = $wb.Worksheets.Item($cfg.Name)

    for ($i = 0; $i -lt $newRows.Count; $i++) {
        $r = $cfg.StartRow + $i
        $sheet.Cells.Item($r, 1).Value = $newRows[$i][0]
        $sheet.Cells.Item($r, 3).Value = $newRows[$i][1]
    }

    $deleteFirst = $cfg.StartRow + $newRows.Count
    if ($deleteFirst -le $cfg.OldLastRow) {
        $sheet.Range("A" + $deleteFirst + ":C" + $cfg.OldLastRow).EntireRow.Delete()
    }
}

Write-Host "edit complete"
